$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.325.48"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").Value = "'1.593.66"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'214.21"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").Value = "'0.494"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'23.99"
$ws.Range("E8").Value = "  +8.38%  "
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'0.0889"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").Value = "'1.822.21"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "'1.592.55"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "'28.355.28"
$ws.Range("D17").Value = "'63.13"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "'227.26"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "'151.66"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'6.59"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.107"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'1.14"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'3.14"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "'1.398.10"
$ws.Range("E34").Value = "  -3.92%  "
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("E36").Value = "  -5.90%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").Value = "'2.55"
$ws.Range("E39").Value = "  +8.83%  "
$ws.Range("D40").Value = "'0.540"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'0.813"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +7.34%  "
$ws.Range("D45").Value = "'0.985"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "'64.38"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'1.732.67"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "'87.46"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "'0.0526"
$ws.Range("E51").Value = "  +0.18%  "
